$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 and 37: LidoDAOToken / TrustWalletToken swap positions, with updated price/volume data
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "  +4.04%  "

# Updated prices and 1h volume percentages for the rest of the list
$ws.Range("D2").Value = "29.540.59"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "1.598.48"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'211.91"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "'26.77"
$ws.Range("E8").Value = "  +5.89%  "
$ws.Range("E9").Value = "  -5.22%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "1.827.43"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "1.602.18"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "29.562.64"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "'0.535"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "'63.66"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").Value = "'239.48"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "'3.98"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "'2.10"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "'154.72"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "'15.33"
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "'6.39"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("D35").Value = "1.432.82"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "'2.31"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "'0.0493"
$ws.Range("E43").Value = "  +6.99%  "
$ws.Range("D44").Value = "'53.27"
$ws.Range("E44").Value = "  +25.34%  "
$ws.Range("D45").Value = "'0.800"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'0.990"
$ws.Range("E47").Value = "  +19.52%  "
$ws.Range("D48").Value = "'65.63"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").Value = "'5.31"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").Value = "1.738.25"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "'86.27"
$ws.Range("E51").Value = "  +1.59%  "
